$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Ideal" starting task count corrected from 41 to 39.
# Downstream formula cells (B4:B30) recompute automatically.
$ws.Range("B3").Value = 39

# The last "Ideal" formula was simplified (no more MAX/1-decimal rounding).
$ws.Range("B31").Formula = "=ROUND(B30-(B$3/28),0)"

# "Actual" values updated for today's work.
$ws.Range("C3:C8").Value = 39
$ws.Range("C12:C27").Value = 22
$ws.Range("C28").Value = 16
$ws.Range("C29").Value = 4

# Move the active selection / cursor to reflect where work left off today.
$ws.Range("B31").Select() | Out-Null
